# Regenerate the handback status report: drop the entries for the
# e464a940-36a9-4bee-a998-9bc4de6759d0 source file (row 3 on every sheet)
# and refresh the processing timestamps recorded for the file that remains
# (row 2, columns E/H) on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": remove row 3 (the e464a940... summary row) and its
# hyperlink, keeping the row 2 hyperlink intact.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows.Item(3).Delete()

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn": remove row 3, keep row 2 hyperlinks, refresh timestamps.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows.Item(3).Delete()

$wsZhCn.Range("E2").Value = "2016-03-11 18:33:42"
$wsZhCn.Range("H2").Value = "2016-03-11 18:34:10"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e08bd87dba2db99d32b6aeee3d484b80778344a2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.zh-cn.xlf", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ad440bb342af6d37a87246565b00caa34c353763/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d609d04d18b7a49195f24d4e556069989c8621c5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.zh-cn.xlf", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de": remove row 3, keep row 2 hyperlinks, refresh timestamps.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Range("E2").Value = "2016-03-11 18:33:45"
$wsDeDe.Range("H2").Value = "2016-03-11 18:34:16"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/1899c4315409e5a9cc6ce1e79cfabb8456889f58/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/487f785d7c0172c05129dc0b4f790bed39accb9e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.de-de.xlf", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/048aa54ae6b221acb9ae568f084b34ebda184e02/e2e/23efbcc1-5d72-4632-9804-a845cd4bc113.md", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dfa18673c9cd7f0d723c813087ae62954778583b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.de-de.xlf", "", "", "23efbcc1-5d72-4632-9804-a845cd4bc113.65693d7e9d6d2756483ea09bc9471a7181b9d112.de-de.xlf")
